$wb = $excel.ActiveWorkbook

$wsCom = $wb.Worksheets.Item("Comunicaciones")
$wsCostos = $wb.Worksheets.Item("Costos")

# --- Comunicaciones (sheet2): fill in symmetric distance matrix, replacing '-' placeholders ---
$wsCom.Range("B2").Value = 0
$wsCom.Range("C2").Value = 2.8
$wsCom.Range("D2").Value = 1.7
$wsCom.Range("E2").Value = 1.2
$wsCom.Range("F2").Value = 0.9
$wsCom.Range("G2").Value = 2.2

$wsCom.Range("B3").Value = 2.8
$wsCom.Range("C3").Value = 0
$wsCom.Range("D3").Value = 1
$wsCom.Range("E3").Value = 0.7
$wsCom.Range("F3").Value = 0.6
$wsCom.Range("G3").Value = 3.1

$wsCom.Range("B4").Value = 1.7
$wsCom.Range("C4").Value = 1
$wsCom.Range("D4").Value = 0
$wsCom.Range("E4").Value = 2.3
$wsCom.Range("F4").Value = 2.5
$wsCom.Range("G4").Value = 1.7

$wsCom.Range("B5").Value = 1.2
$wsCom.Range("C5").Value = 0.7
$wsCom.Range("D5").Value = 2.3
$wsCom.Range("E5").Value = 0
$wsCom.Range("F5").Value = 0.8
$wsCom.Range("G5").Value = 0.9

$wsCom.Range("B6").Value = 0.9
$wsCom.Range("C6").Value = 0.6
$wsCom.Range("D6").Value = 2.5
$wsCom.Range("E6").Value = 0.8
$wsCom.Range("F6").Value = 0
$wsCom.Range("G6").Value = 1.5

# New row 7 - Economía
$wsCom.Range("A7").Value = "Economía"
$wsCom.Range("B7").Value = 2.2
$wsCom.Range("C7").Value = 3.1
$wsCom.Range("D7").Value = 1.7
$wsCom.Range("E7").Value = 0.9
$wsCom.Range("F7").Value = 1.5
$wsCom.Range("G7").Value = 0

# --- Costos (sheet3): fill in symmetric cost matrix, replacing '-' placeholders ---
$wsCostos.Range("B3").Value = 10
$wsCostos.Range("B4").Value = 10
$wsCostos.Range("C4").Value = 11
$wsCostos.Range("B5").Value = 9
$wsCostos.Range("C5").Value = 7
$wsCostos.Range("D5").Value = 8

# --- Selections / active sheet ---
$wsCostos.Range("F11").Select()
$wsCom.Activate()
$wsCom.Range("J15").Select()
